$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for the two new homework columns
$ws.Range("G2").Value = "Homework 2"
$ws.Range("H2").Value = "Homework 3"

# Grades for Homework 2 and Homework 3 per student row
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1

$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1

$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1

$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1

$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1

$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

$ws.Range("H10").Select()
